$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set: year label in column A, value in column B, starting at row 2
$data = @(
    @("2010年", 138904.7),
    @("2011年", 604171.3204),
    @("2013年", 465044.8),
    @("2014年", 154652.8717),
    @("2015年", 115737.3),
    @("2016年", 102800.4),
    @("2017年", 39663.7),
    @("2018年", 39367.3),
    @("2019年", 347100.8385),
    @("2020年", 205726.0892),
    @("2021年", 206049)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused rows 13 through 20
$ws.Range("A13:B20").EntireRow.Delete()
